# Applies the "Importando Modelo Original" update to params.xlsx:
#   - params sheet: refresh Min/Max bounds for the advertising-effectiveness,
#     contact-rate and adoption-fraction rows, and re-apply the sheet's font
#     (as happens when the values are pasted in from the freshly-imported
#     source model) + move the active selection.
#   - levers sheet: collapse the lever table down to levers 1 and 2 (lever 2
#     is repurposed to the "no advertising" baseline) and drop the now-unused
#     levers 3, 4 and the NADV row.

$wb = $excel.ActiveWorkbook

# --- Sheet "params" ---
$params = $wb.Worksheets.Item("params")

# Row 2 - aAdvertisingEffectiveness
$params.Range("C2").Value = 0
$params.Range("D2").Value = 0.02

# Row 3 - aContactRate
$params.Range("C3").Value = 50
$params.Range("D3").Value = 200

# Row 4 - aAdoptionFraction
$params.Range("C4").Value = 0
$params.Range("D4").Value = 0.03

# Re-apply the font across the table (mirrors the font refresh that came
# along with the imported values) and move the active selection.
$params.Range("A1:E5").Font.ThemeColor = 1
$params.Range("B6").Font.ThemeColor = 1
$null = $params.Range("C4").Select()

# --- Sheet "levers" ---
$levers = $wb.Worksheets.Item("levers")

# Lever 2 now points at intensity 2 ...
$levers.Range("D2").Value = 2
# ... and lever "3" becomes the NADV (no-advertising) baseline, intensity 0.
$levers.Range("C3").Value = 0
$levers.Range("D3").Value = 0

# Drop the old levers 3 & 4 rows and the trailing NADV row - the table now
# only spans rows 1-3.
$levers.Rows.Item(4).Delete()
$levers.Rows.Item(4).Delete()
$levers.Rows.Item(4).Delete()

Write-Output "edit applied"
